# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.528.04'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '2.479.48'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.74'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.19'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  -1.07%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.73'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0786'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '2.861.71'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.17'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +8.83%  '
$ws.Range('D16').Value = '2.490.22'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.766'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').Value = '41.523.55'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('E20').Value = '  +2.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.65'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +4.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.32'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.14'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.88'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.81%  '
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.67'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.05'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.39'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.46'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0755'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.35'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('E36').Value = '  -8.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.106'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.42%  '
$ws.Range('E38').Value = '  -3.15%  '
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.12'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.58'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.67%  '
$ws.Range('D44').Value = '1.969.85'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.98'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.51%  '
$ws.Range('D48').Value = '2.719.82'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.60'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.03'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.57'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.81%  '
